$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts B:S -> C:T)
$ws.Columns("B:B").Insert()

# Merge the new B9:B10 cell (empty range) first, then paste formats from
# the existing merged C9:C10 range so the style index gets reused instead
# of the engine synthesizing a new split-border style for the merge.
$ws.Range("B9:B10").Merge() | Out-Null
$ws.Range("C9:C10").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New header label for the inserted column
$ws.Range("B9").Value = "CODIGO"

# Restore the selection to match the authored state
$ws.Range("C6:T6").Select() | Out-Null

Write-Output "done"
